$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.97750833333333
$ws.Range("H2").Value = 65.932525
$ws.Range("I2").Value = 0.5427578249542736
$ws.Range("J2").Value = 0.5427578249542736
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.749486333333333
$ws.Range("N2").Value = 11.248459
$ws.Range("O2").Value = 0.07824568942484071
$ws.Range("P2").Value = 0.07824568942484071
$ws.Range("Q2").Value = 82.40436713655278
$ws.Range("R2").Value = 741.639304228975
$ws.Range("S2").Value = 0.04246846020427415
$ws.Range("T2").Value = 0.04246846020427415

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.97750833333333
$ws.Range("H3").Value = 65.932525
$ws.Range("I3").Value = 0.5427578249542736
$ws.Range("J3").Value = 0.5427578249542736
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.933221666666667
$ws.Range("N3").Value = 8.799665000000001
$ws.Range("O3").Value = 0.06121157170352321
$ws.Range("P3").Value = 0.06121157170352321
$ws.Range("Q3").Value = 64.46490362268055
$ws.Range("R3").Value = 580.184132604125
$ws.Range("S3").Value = 0.03322305951983682
$ws.Range("T3").Value = 0.03322305951983682

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.97750833333333
$ws.Range("H4").Value = 65.932525
$ws.Range("I4").Value = 0.5427578249542736
$ws.Range("J4").Value = 0.5427578249542736
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.50096466666667
$ws.Range("N4").Value = 37.502894
$ws.Range("O4").Value = 0.2608748270724658
$ws.Range("P4").Value = 0.2608748270724658
$ws.Range("Q4").Value = 274.7400551363722
$ws.Range("R4").Value = 2472.66049622735
$ws.Range("S4").Value = 0.1415918537271738
$ws.Range("T4").Value = 0.1415918537271738

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.97750833333333
$ws.Range("H5").Value = 65.932525
$ws.Range("I5").Value = 0.5427578249542736
$ws.Range("J5").Value = 0.5427578249542736
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.73572533333333
$ws.Range("N5").Value = 86.207176
$ws.Range("O5").Value = 0.5996679117991702
$ws.Range("P5").Value = 0.5996679117991702
$ws.Range("Q5").Value = 631.5396429777111
$ws.Range("R5").Value = 5683.8567867994
$ws.Range("S5").Value = 0.3254744515029888
$ws.Range("T5").Value = 0.3254744515029888

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.688376
$ws.Range("H6").Value = 38.065128
$ws.Range("I6").Value = 0.3133528721960219
$ws.Range("J6").Value = 0.3133528721960219
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.749486333333333
$ws.Range("N6").Value = 11.248459
$ws.Range("O6").Value = 0.07824568942484071
$ws.Range("P6").Value = 0.07824568942484071
$ws.Range("Q6").Value = 47.57489240419467
$ws.Range("R6").Value = 428.174031637752
$ws.Range("S6").Value = 0.02451851151823174
$ws.Range("T6").Value = 0.02451851151823174

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.688376
$ws.Range("H7").Value = 38.065128
$ws.Range("I7").Value = 0.3133528721960219
$ws.Range("J7").Value = 0.3133528721960219
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.933221666666667
$ws.Range("N7").Value = 8.799665000000001
$ws.Range("O7").Value = 0.06121157170352321
$ws.Range("P7").Value = 0.06121157170352321
$ws.Range("Q7").Value = 37.21781939801333
$ws.Range("R7").Value = 334.96037458212
$ws.Range("S7").Value = 0.01918082180493174
$ws.Range("T7").Value = 0.01918082180493174

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.688376
$ws.Range("H8").Value = 38.065128
$ws.Range("I8").Value = 0.3133528721960219
$ws.Range("J8").Value = 0.3133528721960219
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.50096466666667
$ws.Range("N8").Value = 37.502894
$ws.Range("O8").Value = 0.2608748270724658
$ws.Range("P8").Value = 0.2608748270724658
$ws.Range("Q8").Value = 158.6169400533813
$ws.Range("R8").Value = 1427.552460480432
$ws.Range("S8").Value = 0.08174587634679771
$ws.Range("T8").Value = 0.08174587634679771

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.688376
$ws.Range("H9").Value = 38.065128
$ws.Range("I9").Value = 0.3133528721960219
$ws.Range("J9").Value = 0.3133528721960219
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 28.73572533333333
$ws.Range("N9").Value = 86.207176
$ws.Range("O9").Value = 0.5996679117991702
$ws.Range("P9").Value = 0.5996679117991702
$ws.Range("Q9").Value = 364.6096876620587
$ws.Range("R9").Value = 3281.487188958528
$ws.Range("S9").Value = 0.1879076625260607
$ws.Range("T9").Value = 0.1879076625260607

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.008189666666667
$ws.Range("H10").Value = 6.024569
$ws.Range("I10").Value = 0.04959436889042158
$ws.Range("J10").Value = 0.04959436889042158
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.749486333333333
$ws.Range("N10").Value = 11.248459
$ws.Range("O10").Value = 0.07824568942484071
$ws.Range("P10").Value = 0.07824568942484071
$ws.Range("Q10").Value = 7.529679709907889
$ws.Range("R10").Value = 67.76711738917099
$ws.Range("S10").Value = 0.003880545585420909
$ws.Range("T10").Value = 0.003880545585420909

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.008189666666667
$ws.Range("H11").Value = 6.024569
$ws.Range("I11").Value = 0.04959436889042158
$ws.Range("J11").Value = 0.04959436889042158
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.933221666666667
$ws.Range("N11").Value = 8.799665000000001
$ws.Range("O11").Value = 0.06121157170352321
$ws.Range("P11").Value = 0.06121157170352321
$ws.Range("Q11").Value = 5.890465441042778
$ws.Range("R11").Value = 53.014188969385
$ws.Range("S11").Value = 0.003035749267427021
$ws.Range("T11").Value = 0.003035749267427021

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.008189666666667
$ws.Range("H12").Value = 6.024569
$ws.Range("I12").Value = 0.04959436889042158
$ws.Range("J12").Value = 0.04959436889042158
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.50096466666667
$ws.Range("N12").Value = 37.502894
$ws.Range("O12").Value = 0.2608748270724658
$ws.Range("P12").Value = 0.2608748270724658
$ws.Range("Q12").Value = 25.10430806696511
$ws.Range("R12").Value = 225.938772602686
$ws.Range("S12").Value = 0.01293792240805681
$ws.Range("T12").Value = 0.01293792240805681

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.008189666666667
$ws.Range("H13").Value = 6.024569
$ws.Range("I13").Value = 0.04959436889042158
$ws.Range("J13").Value = 0.04959436889042158
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 28.73572533333333
$ws.Range("N13").Value = 86.207176
$ws.Range("O13").Value = 0.5996679117991702
$ws.Range("P13").Value = 0.5996679117991702
$ws.Range("Q13").Value = 57.70678667857155
$ws.Range("R13").Value = 519.3610801071439
$ws.Range("S13").Value = 0.02974015162951684
$ws.Range("T13").Value = 0.02974015162951684

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.818218
$ws.Range("H14").Value = 11.454654
$ws.Range("I14").Value = 0.09429493395928291
$ws.Range("J14").Value = 0.09429493395928291
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.749486333333333
$ws.Range("N14").Value = 11.248459
$ws.Range("O14").Value = 0.07824568942484071
$ws.Range("P14").Value = 0.07824568942484071
$ws.Range("Q14").Value = 14.31635620868733
$ws.Range("R14").Value = 128.847205878186
$ws.Range("S14").Value = 0.007378172116913917
$ws.Range("T14").Value = 0.007378172116913917

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.818218
$ws.Range("H15").Value = 11.454654
$ws.Range("I15").Value = 0.09429493395928291
$ws.Range("J15").Value = 0.09429493395928291
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.933221666666667
$ws.Range("N15").Value = 8.799665000000001
$ws.Range("O15").Value = 0.06121157170352321
$ws.Range("P15").Value = 0.06121157170352321
$ws.Range("Q15").Value = 11.19967976565667
$ws.Range("R15").Value = 100.79711789091
$ws.Range("S15").Value = 0.005771941111327632
$ws.Range("T15").Value = 0.005771941111327632

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.818218
$ws.Range("H16").Value = 11.454654
$ws.Range("I16").Value = 0.09429493395928291
$ws.Range("J16").Value = 0.09429493395928291
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.50096466666667
$ws.Range("N16").Value = 37.502894
$ws.Range("O16").Value = 0.2608748270724658
$ws.Range("P16").Value = 0.2608748270724658
$ws.Range("Q16").Value = 47.73140830763067
$ws.Range("R16").Value = 429.582674768676
$ws.Range("S16").Value = 0.02459917459043752
$ws.Range("T16").Value = 0.02459917459043752

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.818218
$ws.Range("H17").Value = 11.454654
$ws.Range("I17").Value = 0.09429493395928291
$ws.Range("J17").Value = 0.09429493395928291
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 28.73572533333333
$ws.Range("N17").Value = 86.207176
$ws.Range("O17").Value = 0.5996679117991702
$ws.Range("P17").Value = 0.5996679117991702
$ws.Range("Q17").Value = 109.7192637107893
$ws.Range("R17").Value = 987.4733733971041
$ws.Range("S17").Value = 0.05654564614060385
$ws.Range("T17").Value = 0.05654564614060385

